$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(89, 1).Value = 44322
$ws.Cells.Item(89, 2).Value = 1045
$ws.Cells.Item(89, 3).Value = 2478
$ws.Cells.Item(89, 4).Value = 615
$ws.Cells.Item(89, 5).Value = 5425
$ws.Cells.Item(89, 6).Value = 189
$ws.Cells.Item(89, 7).Value = 1929
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 28500
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 3638
$ws.Cells.Item(89, 14).Value = 310
$ws.Cells.Item(89, 15).Value = 925
$ws.Cells.Item(89, 16).Value = 41971
$ws.Cells.Item(90, 1).Value = 44323
$ws.Cells.Item(90, 2).Value = 450
$ws.Cells.Item(90, 3).Value = 2888
$ws.Cells.Item(90, 4).Value = 850
$ws.Cells.Item(90, 5).Value = 5325
$ws.Cells.Item(90, 6).Value = 724
$ws.Cells.Item(90, 7).Value = 1928
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 28500
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = 3638
$ws.Cells.Item(90, 14).Value = 11
$ws.Cells.Item(90, 15).Value = 925
$ws.Cells.Item(90, 16).Value = 42279
$ws.Cells.Item(91, 1).Value = 44326
$ws.Cells.Item(91, 2).Value = 671
$ws.Cells.Item(91, 3).Value = 3126
$ws.Cells.Item(91, 4).Value = 850
$ws.Cells.Item(91, 5).Value = 5275
$ws.Cells.Item(91, 6).Value = 686
$ws.Cells.Item(91, 7).Value = 1931
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 28500
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = 3638
$ws.Cells.Item(91, 14).Value = 135
$ws.Cells.Item(91, 15).Value = 878
$ws.Cells.Item(91, 16).Value = 42470
$ws.Cells.Item(92, 1).Value = 44327
$ws.Cells.Item(92, 2).Value = 530
$ws.Cells.Item(92, 3).Value = 2946
$ws.Cells.Item(92, 4).Value = 800
$ws.Cells.Item(92, 5).Value = 5275
$ws.Cells.Item(92, 6).Value = 208
$ws.Cells.Item(92, 7).Value = 1938
$ws.Cells.Item(92, 8).Value = 1000
$ws.Cells.Item(92, 9).Value = 28500
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 3638
$ws.Cells.Item(92, 14).Value = 224
$ws.Cells.Item(92, 15).Value = 837
$ws.Cells.Item(92, 16).Value = 42298
$ws.Cells.Item(93, 1).Value = 44328
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(93, 3).Value = 2946
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 5275
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 1938
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 28500
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 3638
$ws.Cells.Item(93, 14).Value = 0
$ws.Cells.Item(93, 15).Value = 837
$ws.Cells.Item(93, 16).Value = 42298
$ws.Cells.Item(94, 1).Value = 44333
$ws.Cells.Item(94, 2).Value = 741
$ws.Cells.Item(94, 3).Value = 3016
$ws.Cells.Item(94, 4).Value = 1172
$ws.Cells.Item(94, 5).Value = 5677
$ws.Cells.Item(94, 6).Value = 781
$ws.Cells.Item(94, 7).Value = 2034
$ws.Cells.Item(94, 8).Value = 1000
$ws.Cells.Item(94, 9).Value = 28500
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 1235
$ws.Cells.Item(94, 13).Value = 3674
$ws.Cells.Item(94, 14).Value = 436
$ws.Cells.Item(94, 15).Value = 1138
$ws.Cells.Item(94, 16).Value = 42901
